# The commit re-orders the 21 observation rows (rows 2-22) of the sheet:
# each row's full set of column values moves to a different row position;
# no rows are added or removed and the header row (row 1) is untouched.
#
# Mapping below: key = destination (new) sheet row, value = source (old)
# sheet row whose data should end up there. Derived from the authoritative
# XML diff (every column for every row matches this permutation exactly).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowMap = @{
    2  = 16
    3  = 18
    4  = 20
    5  = 13
    6  = 8
    7  = 22
    8  = 12
    9  = 17
    10 = 5
    11 = 19
    12 = 9
    13 = 14
    14 = 4
    15 = 21
    16 = 7
    17 = 2
    18 = 15
    19 = 3
    20 = 11
    21 = 6
    22 = 10
}

$firstRow = 2
$lastRow = 22
$firstCol = 1   # A
$lastCol = 51   # AY

# Snapshot the current ("before") values of every data row/column so the
# permutation can be computed purely in memory before anything is written
# back (avoids clobbering source rows we still need to read).
$srcRange = $ws.Range($ws.Cells.Item($firstRow, $firstCol), $ws.Cells.Item($lastRow, $lastCol))
$srcValues = $srcRange.Value()

$numRows = $lastRow - $firstRow + 1
$numCols = $lastCol - $firstCol + 1

# Columns Y and AA hold plain-text dates (e.g. "2023-08-31"). They are
# stored as literal text in the source file, but Excel's COM layer
# auto-converts date-shaped text to a real date serial the moment it is
# assigned. Pre-mark those two columns with the classic leading-apostrophe
# "treat as text" prefix while building the array, so the single bulk
# write below lands as text (the prefix itself is stripped by Excel and
# never appears in the stored value).
$dateCols = @(25, 27)  # Y, AA (1-based column numbers)

# NOTE: a freshly allocated .NET array (New-Object 'object[,]') is
# zero-based, while the array handed back by a COM Range's .Value() is
# one-based (row/col 1 is the first cell) - offset accordingly below.
$newValues = New-Object 'object[,]' $numRows, $numCols

foreach ($destRow in $rowMap.Keys) {
    $sourceRow = $rowMap[$destRow]
    $destIdx = $destRow - $firstRow      # 0-based for $newValues
    $sourceIdx = $sourceRow - $firstRow + 1  # 1-based for $srcValues
    for ($c = 1; $c -le $numCols; $c++) {
        $val = $srcValues[$sourceIdx, $c]
        $absoluteCol = $firstCol + $c - 1
        if (($dateCols -contains $absoluteCol) -and ($val -ne $null) -and ($val -is [string]) -and ($val.Length -gt 0)) {
            $val = "'" + $val
        }
        $newValues[$destIdx, ($c - 1)] = $val
    }
}

$srcRange.Value = $newValues
